$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CO2_Fossil_Fuel")
$ws.Activate()

# Move the source note and units note from column F to column I
$ws.Range("I1").Value2 = $ws.Range("F1").Value2
$ws.Range("F1").ClearContents()

$ws.Range("I2").Value2 = $ws.Range("F2").Value2
$ws.Range("F2").ClearContents()

# Update the active selection to match the new view state
$ws.Range("P30").Select()
